$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(26).Insert()

$ws.Cells.Item(26, 1).Value = 9
$ws.Cells.Item(26, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(26, 3).Value = 'Metropolitana'
$ws.Cells.Item(26, 4).Value = 44725
$ws.Cells.Item(26, 5).Value = 13
$ws.Cells.Item(26, 6).Value = 100114007
$ws.Cells.Item(26, 7).Value = 'Jengibre'
$ws.Cells.Item(26, 8).Value = 'Sin especificar'
$ws.Cells.Item(26, 9).Value = 'Primera'
$ws.Cells.Item(26, 10).Value = 610
$ws.Cells.Item(26, 11).Value = 14000
$ws.Cells.Item(26, 12).Value = 15000
$ws.Cells.Item(26, 13).Value = 14500
$ws.Cells.Item(26, 14).Value = '$/caja 13 kilos'
$ws.Cells.Item(26, 15).Value = 'Perú'
$ws.Cells.Item(26, 16).Value = 1115
$ws.Cells.Item(26, 17).Value = 13
$ws.Cells.Item(26, 18).Value = 'Hortaliza'
